$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 310
$ws.Range("F5").Value = 1258
$ws.Range("F7").Value = 299
$ws.Range("F8").Value = 1097
$ws.Range("F9").Value = 427
$ws.Range("F10").Value = 6895
$ws.Range("F11").Value = 79
$ws.Range("F14").Value = 7800
$ws.Range("F17").Value = 5124
$ws.Range("F20").Value = 967
$ws.Range("F21").Value = 4538
$ws.Range("F23").Value = 366
$ws.Range("F26").Value = 283
$ws.Range("F29").Value = 2014
$ws.Range("F35").Value = 1373
$ws.Range("F36").Value = 20
$ws.Range("F37").Value = 2098

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 28
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 23

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 251

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 251
$ws.Range("F7").Value = 310
$ws.Range("F8").Value = 1258
$ws.Range("F11").Value = 299
$ws.Range("F12").Value = 1097
$ws.Range("F13").Value = 427
$ws.Range("F14").Value = 6895
$ws.Range("F15").Value = 79
$ws.Range("F18").Value = 7800
$ws.Range("F21").Value = 5124
$ws.Range("F24").Value = 967
$ws.Range("F25").Value = 4538
$ws.Range("F27").Value = 366
$ws.Range("F31").Value = 28
$ws.Range("F32").Value = 283
$ws.Range("F35").Value = 2014
$ws.Range("F41").Value = 16
$ws.Range("F42").Value = 1373
$ws.Range("F43").Value = 20
$ws.Range("F44").Value = 2098
$ws.Range("F45").Value = 23
